$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (@fresheroffcampusdiscussion) updates
$ws.Range("E11").Value = "2026-02-22T13:50:31.876368+00:00"
$ws.Range("H11").Value = 7
$ws.Range("I11").Value = 3
$ws.Range("L11").Value = "[68050, 68073, 68060, 68065, 68192, 68188, 68209]"
$ws.Range("M11").Value = "[68188, 68187, 68199]"
